$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update existing data in Sheet1
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 5
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = 6
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = 7
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = 8

# Add a new empty sheet named Sheet2 right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Restore Sheet1 as the active/selected sheet
$ws1.Select()
